$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 500657
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 500657
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 500657
$ws.Range("N3").Value = -500885
# Row 39
$ws.Range("H39").Value = 159.36363
$ws.Range("I39").Value = 69
$ws.Range("J39").Value = 317.5
$ws.Range("K39").Value = 207
$ws.Range("L39").Value = 952.5
$ws.Range("M39").Value = 89
# Row 40
$ws.Range("H40").Value = 2041.8485
$ws.Range("I40").Value = 1951.8572
$ws.Range("J40").Value = 2199.3333
$ws.Range("K40").Value = 1951.8572
$ws.Range("L40").Value = 2199.3333
$ws.Range("M40").Value = -1776.8572
$ws.Range("N40").Value = -2549.3333
# Row 42
$ws.Range("H42").Value = 274.14285
$ws.Range("I42").Value = 45.714287
$ws.Range("J42").Value = 502.57144
$ws.Range("K42").Value = 137.142861
$ws.Range("L42").Value = 1507.71432
$ws.Range("M42").Value = 92.85713900000002
$ws.Range("N42").Value = -1967.71432
# Row 52
$ws.Range("H52").Value = 2710
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 2710
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 8130
$ws.Range("N52").Value = -8450
# Row 101
$ws.Range("H101").Value = 513.3333
$ws.Range("I101").Value = 416
$ws.Range("J101").Value = 1000
$ws.Range("K101").Value = 1248
$ws.Range("L101").Value = 3000
$ws.Range("M101").Value = 374
# Row 102
$ws.Range("H102").Value = 500657
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 500657
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 500657
$ws.Range("N102").Value = -507147
# Row 131
$ws.Range("H131").Value = 5858.154
$ws.Range("I131").Value = 1315.6
$ws.Range("J131").Value = 21000
$ws.Range("K131").Value = 3946.8
$ws.Range("L131").Value = 63000
$ws.Range("M131").Value = 1093.2
$ws.Range("N131").Value = -73080
# Row 132
$ws.Range("H132").Value = 410890.4
$ws.Range("I132").Value = 553390.75
$ws.Range("J132").Value = 62556.223
$ws.Range("K132").Value = 1660172.25
$ws.Range("L132").Value = 187668.669
$ws.Range("M132").Value = -1657642.25
$ws.Range("N132").Value = -192728.669

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 132
$ws.Range("H132").Value = 3073.875
$ws.Range("I132").Value = 2733.4
$ws.Range("J132").Value = 4289.857
$ws.Range("K132").Value = 8200.200000000001
$ws.Range("L132").Value = 12869.571
$ws.Range("M132").Value = -5670.200000000001
$ws.Range("N132").Value = -17929.571
# Row 141
$ws.Range("H141").Value = 30750
$ws.Range("I141").Value = 30000
$ws.Range("J141").Value = 31000
$ws.Range("K141").Value = 30000
$ws.Range("L141").Value = 31000
$ws.Range("M141").Value = -24820
$ws.Range("N141").Value = -41360

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 59
$ws.Range("H59").Value = 47825
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 47825
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 47825
$ws.Range("N59").Value = -49519
# Row 61
$ws.Range("H61").Value = 8333.333000000001
$ws.Range("I61").Value = 7000
$ws.Range("J61").Value = 8500
$ws.Range("K61").Value = 7000
$ws.Range("L61").Value = 8500
$ws.Range("M61").Value = -6687
$ws.Range("N61").Value = -9126
# Row 75
$ws.Range("H75").Value = 75289.46000000001
$ws.Range("I75").Value = 5391.8887
$ws.Range("J75").Value = 232559
$ws.Range("K75").Value = 5391.8887
$ws.Range("L75").Value = 232559
$ws.Range("M75").Value = -4455.8887
# Row 78
$ws.Range("H78").Value = 75289.46000000001
$ws.Range("I78").Value = 5391.8887
$ws.Range("J78").Value = 232559
$ws.Range("K78").Value = 16175.6661
$ws.Range("L78").Value = 697677
$ws.Range("M78").Value = -11495.6661
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1146.721
$ws.Range("I31").Value = 892.6053000000001
$ws.Range("J31").Value = 3078
$ws.Range("K31").Value = 892.6053000000001
$ws.Range("L31").Value = 3078
$ws.Range("M31").Value = -597.6053000000001
$ws.Range("N31").Value = -3668
# Row 33
$ws.Range("H33").Value = 11888
$ws.Range("I33").Value = 11888
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 11888
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -11509
# Row 34
$ws.Range("H34").Value = 1146.721
$ws.Range("I34").Value = 892.6053000000001
$ws.Range("J34").Value = 3078
$ws.Range("K34").Value = 892.6053000000001
$ws.Range("L34").Value = 3078
$ws.Range("M34").Value = -690.6053000000001
$ws.Range("N34").Value = -3482
# Row 51
$ws.Range("H51").Value = 29400
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 29400
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 29400
$ws.Range("N51").Value = -30872
# Row 61
$ws.Range("H61").Value = 29400
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 29400
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 29400
$ws.Range("N61").Value = -30096

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 501.8125
$ws.Range("I107").Value = 441.25
$ws.Range("J107").Value = 562.375
$ws.Range("K107").Value = 1323.75
$ws.Range("L107").Value = 1687.125
$ws.Range("M107").Value = 596.25
$ws.Range("N107").Value = -5527.125
# Row 127
$ws.Range("H127").Value = 2457.1428
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 2457.1428
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 7371.428400000001
$ws.Range("N127").Value = -17291.4284
# Row 133
$ws.Range("H133").Value = 9250
$ws.Range("I133").Value = 1000
$ws.Range("J133").Value = 12000
$ws.Range("K133").Value = 3000
$ws.Range("L133").Value = 36000
$ws.Range("M133").Value = 2060
$ws.Range("N133").Value = -46120
# Row 137
$ws.Range("H137").Value = 4813222.5
$ws.Range("I137").Value = 10002109
$ws.Range("J137").Value = 96052.37
$ws.Range("K137").Value = 30006327
$ws.Range("L137").Value = 288157.11
$ws.Range("M137").Value = -30001227

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 139
$ws.Range("H139").Value = 34542
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 34542
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 34542
$ws.Range("N139").Value = -44822

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1216.8334
$ws.Range("I46").Value = 1157.1428
$ws.Range("J46").Value = 1300.4
$ws.Range("K46").Value = 1157.1428
$ws.Range("L46").Value = 1300.4
$ws.Range("M46").Value = -969.1428000000001
$ws.Range("N46").Value = -1676.4
# Row 55
$ws.Range("H55").Value = 403.3125
$ws.Range("I55").Value = 322.125
$ws.Range("J55").Value = 484.5
$ws.Range("K55").Value = 322.125
$ws.Range("L55").Value = 484.5
$ws.Range("M55").Value = -149.125
$ws.Range("N55").Value = -830.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 9
$ws.Range("H9").Value = 43000
$ws.Range("I9").Value = 43000
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 43000
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -42860
$ws.Range("N9").ClearContents()
# Row 138
$ws.Range("H138").Value = 69999.5
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 69999.5
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 69999.5
$ws.Range("N138").Value = -80279.5
